# Apply cryptos list update (price/volume refresh + row reorder for rows 25-26 and 46-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.020.66"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "'2.661.75"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'580.20"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "'145.71"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'6.59"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "'0.382"
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "'3.126.51"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").Value = "'25.78"
$ws.Range("E14").Value = "  +10.64%  "
$ws.Range("D15").Value = "'60.913.99"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "'0.0000144"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "'2.660.37"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "'11.61"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "'4.74"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "'351.67"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'6.99"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "'0.537"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "'64.13"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.162"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'8.21"
$ws.Range("E27").Value = "  +6.20%  "
$ws.Range("E28").Value = "  +8.88%  "
$ws.Range("D29").Value = "'0.0₃0817"
$ws.Range("E29").Value = "  +3.74%  "
$ws.Range("D30").Value = "'6.77"
$ws.Range("E30").Value = "  +5.90%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'167.04"
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("D33").Value = "'19.95"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").Value = "'1.07"
$ws.Range("E34").Value = "  +9.67%  "
$ws.Range("D35").Value = "'4.47"
$ws.Range("E35").Value = "  +5.98%  "
$ws.Range("D36").Value = "'1.32"
$ws.Range("E36").Value = "  +8.38%  "
$ws.Range("D37").Value = "'1.65"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").Value = "'329.05"
$ws.Range("E38").Value = "  +11.21%  "
$ws.Range("D39").Value = "'4.02"
$ws.Range("E39").Value = "  +5.06%  "
$ws.Range("D40").Value = "'38.44"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").Value = "'0.883"
$ws.Range("E41").Value = "  +3.96%  "
$ws.Range("E42").Value = "  +8.33%  "
$ws.Range("D43").Value = "'20.56"
$ws.Range("E43").Value = "  +4.26%  "
$ws.Range("D44").Value = "'135.26"
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("D45").Value = "'0.100"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0562"
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.616"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'20.60"
$ws.Range("E48").Value = "  +3.79%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0247"
$ws.Range("E49").Value = "  +3.00%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'2.147.64"
$ws.Range("E50").Value = "  +5.96%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'0.996"
$ws.Range("E51").Value = "  -0.18%  "
